$wb = $excel.ActiveWorkbook

# --- Intro sheet: text update, new B7 styled cell, dimension/selection changes ---
$wsIntro = $wb.Worksheets.Item("Intro")
$wsIntro.Activate()
$wsIntro.Range("A2").Value2 = 'tn:The Energy tab displays a roast''s energy consumption.   CO2 emissions are also calculated to monitor the impact of the roasting operation.  Settings must be made for each energy load.  Loads are the main burners, motors and blowers, and an afterburner if one is used.  The energy used for pre-heating, between batch, and roaster cooling protocols are included in the calculations, and settings are available for them as well.\n\nNote that pre-heating and roaster cooling energy values are applied to the first roast of a roasting session.  Between batch energies are applied to every roast except the first.  Tick the "Between batches after Pre-Heating box to apply the between batch value to the first roast.\n\nFollow the steps below to set the energy inputs for the roast machine and afterburner.'
$wsIntro.Range("A7").Copy()
$wsIntro.Range("B7").PasteSpecial(-4122)
$wsIntro.Range("A2").Select()

# --- Details sheet: selection moves to B24 ---
$wsDetails = $wb.Worksheets.Item("Details")
$wsDetails.Activate()
$wsDetails.Range("B24").Select()

# --- Loads sheet: no content or view changes needed ---

# --- Protocol sheet: text updates, becomes the active/selected tab ---
$wsProtocol = $wb.Worksheets.Item("Protocol")
$wsProtocol.Activate()
$wsProtocol.Range("A2").Value2 = 'tn:The Protocol settings allow including Pre-Heating, Between Batch (BBP) and Cooling protocol energy consumption.  There are two ways to specify these values.  The first assumes a constant load setting for a defined period of time.  An example for pre-heating is to set a Duration of 45:00 (45 minutes) at 30% Burner setting.  Percentages must be entered with the percent sign (30%).  When a percentage is entered  a corresponding Duration must be entered.'
$wsProtocol.Range("B16").Value2 = 'This box should be ticked when a Between Batches protocol run is done after the Pre-heating and before the roast.'
$wsProtocol.Range("A2").Select()
